# Auto-generated Excel COM-interop script
# Applies numeric updates to ALC, ARM, BSM, CRP, CUL, GSM, WVR sheets
# per scheduled market-price refresh (see commit message).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3860.6
$ws.Range("H77").Value = 3860.6
$ws.Range("H112").Value = 1255.9791
$ws.Range("J112").Value = 1318.7675
$ws.Range("L112").Value = 3956.3025
$ws.Range("N112").Value = -6172.3025
$ws.Range("H125").Value = 1071342.9
$ws.Range("I125").Value = 2114
$ws.Range("J125").Value = 1962367
$ws.Range("K125").Value = 19026
$ws.Range("L125").Value = 17661303
$ws.Range("M125").Value = -16566
$ws.Range("N125").Value = -17666223
$ws.Range("H129").Value = 814
$ws.Range("J129").Value = 897.6923
$ws.Range("L129").Value = 2693.0769
$ws.Range("N129").Value = -12693.0769

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4592.278
$ws.Range("I74").Value = 5740.36
$ws.Range("J74").Value = 1983
$ws.Range("K74").Value = 5740.36
$ws.Range("L74").Value = 1983
$ws.Range("M74").Value = -4866.36
$ws.Range("N74").Value = -3731
$ws.Range("H77").Value = 4592.278
$ws.Range("I77").Value = 5740.36
$ws.Range("J77").Value = 1983
$ws.Range("K77").Value = 28701.8
$ws.Range("L77").Value = 9915
$ws.Range("M77").Value = -24333.8
$ws.Range("N77").Value = -18651

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("N94").ClearContents()
$ws.Range("H97").Value = 950
$ws.Range("I97").Value = 950
$ws.Range("K97").Value = 950
$ws.Range("M97").Value = 41
$ws.Range("H107").Value = 1066.6666
$ws.Range("I107").Value = 700
$ws.Range("K107").Value = 700
$ws.Range("M107").Value = 1220
$ws.Range("H134").Value = 2014.5714
$ws.Range("I134").Value = 1593.2858
$ws.Range("J134").Value = 2857.1428
$ws.Range("K134").Value = 4779.857400000001
$ws.Range("L134").Value = 8571.428400000001
$ws.Range("M134").Value = -2244.857400000001
$ws.Range("N134").Value = -13641.4284
$ws.Range("H135").Value = 38966.668
$ws.Range("J135").Value = 38966.668
$ws.Range("L135").Value = 38966.668
$ws.Range("N135").Value = -49106.668
$ws.Range("H140").Value = 41158.95
$ws.Range("J140").Value = 41158.95
$ws.Range("L140").Value = 41158.95
$ws.Range("N140").Value = -51518.95

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4073.116
$ws.Range("I31").Value = 2298.913
$ws.Range("J31").Value = 4960.2173
$ws.Range("K31").Value = 2298.913
$ws.Range("L31").Value = 4960.2173
$ws.Range("M31").Value = -2003.913
$ws.Range("N31").Value = -5550.2173
$ws.Range("H34").Value = 4073.116
$ws.Range("I34").Value = 2298.913
$ws.Range("J34").Value = 4960.2173
$ws.Range("K34").Value = 2298.913
$ws.Range("L34").Value = 4960.2173
$ws.Range("M34").Value = -2096.913
$ws.Range("N34").Value = -5364.2173
$ws.Range("H62").Value = 11001.667
$ws.Range("I62").Value = 10005
$ws.Range("J62").Value = 11500
$ws.Range("K62").Value = 10005
$ws.Range("L62").Value = 11500
$ws.Range("M62").Value = -9381
$ws.Range("N62").Value = -12748
$ws.Range("H65").Value = 11001.667
$ws.Range("I65").Value = 10005
$ws.Range("J65").Value = 11500
$ws.Range("K65").Value = 50025
$ws.Range("L65").Value = 57500
$ws.Range("M65").Value = -46905
$ws.Range("N65").Value = -63740
$ws.Range("H99").Value = 7929.3335
$ws.Range("I99").Value = 1380
$ws.Range("J99").Value = 50500
$ws.Range("K99").Value = 1380
$ws.Range("L99").Value = 50500
$ws.Range("M99").Value = 118
$ws.Range("N99").Value = -53496
$ws.Range("H126").Value = 7929.3335
$ws.Range("I126").Value = 1380
$ws.Range("J126").Value = 50500
$ws.Range("K126").Value = 4140
$ws.Range("L126").Value = 151500
$ws.Range("M126").Value = -1670
$ws.Range("N126").Value = -156440
$ws.Range("H134").Value = 2070.4333
$ws.Range("I134").Value = 2297
$ws.Range("J134").Value = 1447.375
$ws.Range("K134").Value = 6891
$ws.Range("L134").Value = 4342.125
$ws.Range("M134").Value = -4356
$ws.Range("N134").Value = -9412.125
$ws.Range("H140").Value = 56564.832
$ws.Range("J140").Value = 56564.832
$ws.Range("L140").Value = 56564.832
$ws.Range("N140").Value = -66924.83199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1596425.1
$ws.Range("I5").Value = 702
$ws.Range("K5").Value = 2106
$ws.Range("M5").Value = -1994
$ws.Range("H11").Value = 138.16667
$ws.Range("I11").Value = 90
$ws.Range("J11").Value = 186.33333
$ws.Range("K11").Value = 270
$ws.Range("L11").Value = 558.99999
$ws.Range("M11").Value = -130
$ws.Range("N11").Value = -838.99999
$ws.Range("H87").Value = 12450
$ws.Range("I87").Value = 12450
$ws.Range("K87").Value = 37350
$ws.Range("M87").Value = -36102
$ws.Range("H90").Value = 12450
$ws.Range("I90").Value = 12450
$ws.Range("K90").Value = 112050
$ws.Range("M90").Value = -105810
$ws.Range("H107").Value = 795373.5600000001
$ws.Range("I107").Value = 382
$ws.Range("J107").Value = 1126620.1
$ws.Range("K107").Value = 1146
$ws.Range("L107").Value = 3379860.3
$ws.Range("M107").Value = 774
$ws.Range("N107").Value = -3383700.3
$ws.Range("H131").Value = 852.52
$ws.Range("J131").Value = 881.2766
$ws.Range("L131").Value = 2643.8298
$ws.Range("N131").Value = -12723.8298
$ws.Range("H135").Value = 1596425.1
$ws.Range("I135").Value = 702
$ws.Range("K135").Value = 6318
$ws.Range("M135").Value = -3783

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H138").Value = 42218.445
$ws.Range("J138").Value = 42218.445
$ws.Range("L138").Value = 42218.445
$ws.Range("N138").Value = -52498.445
$ws.Range("H140").Value = 39043916
$ws.Range("J140").Value = 39043916
$ws.Range("L140").Value = 39043916
$ws.Range("N140").Value = -39054276

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1138.6
$ws.Range("I113").Value = 1814.3334
$ws.Range("J113").Value = 125
$ws.Range("K113").Value = 5443.0002
$ws.Range("L113").Value = 375
$ws.Range("M113").Value = -3273.0002
$ws.Range("N113").Value = -4715
$ws.Range("H132").Value = 2386.9167
$ws.Range("I132").Value = 1610.2307
$ws.Range("J132").Value = 3304.818
$ws.Range("K132").Value = 4830.6921
$ws.Range("L132").Value = 9914.454000000002
$ws.Range("M132").Value = -2300.6921
$ws.Range("N132").Value = -14974.454
$ws.Range("H137").Value = 71422
$ws.Range("J137").Value = 71422
$ws.Range("L137").Value = 71422
$ws.Range("N137").Value = -81622
